$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows: Ethereum after Gold (new row 22), Tehther after Mexico (new row 26) ---
$ws.Range("A22").EntireRow.Insert()
$ws.Range("A26").EntireRow.Insert()

# --- Fix "NG" label -> "NaturalGas" ---
$ws.Range("A4").Value = "NaturalGas"

# --- Set labels for newly inserted rows ---
$ws.Range("A22").Value = "Ethereum"
$ws.Range("A26").Value = "Tehther"

# --- Update numeric data cells (B:H) for rows 2-28 with the new descriptive statistics ---
# row 2
$ws.Range("B2").Value = 0.00048304055142385797
$ws.Range("C2").Value = 0.020422206523585951
$ws.Range("D2").Value = 0.001360675502368736
$ws.Range("E2").Value = -0.15993026567604399
$ws.Range("F2").Value = 0.13022281024270829
$ws.Range("G2").Value = -0.01036455869663208
$ws.Range("H2").Value = 1.257260545190098

# row 3
$ws.Range("B3").Value = -0.00031116379063839459
$ws.Range("C3").Value = 0.016469670058470439
$ws.Range("D3").Value = 0.0001067841228143962
$ws.Range("E3").Value = -0.12854533819572289
$ws.Range("F3").Value = 0.088936485781973218
$ws.Range("G3").Value = 0.022959496288078252
$ws.Range("H3").Value = 1.5292908219206081

# row 4
$ws.Range("B4").Value = 0.0010804943372392599
$ws.Range("C4").Value = 0.045277351461191447
$ws.Range("D4").Value = 0.00034240270918217691
$ws.Range("E4").Value = -0.30047976257842501
$ws.Range("F4").Value = 0.4355212407018601
$ws.Range("G4").Value = 0.035583476779334607
$ws.Range("H4").Value = 1.3750363757530599

# row 5
$ws.Range("B5").Value = 0.00073111956776918759
$ws.Range("C5").Value = 0.016191306748179441
$ws.Range("D5").Value = 0.0011760739519850461
$ws.Range("E5").Value = -0.096032648614778005
$ws.Range("F5").Value = 0.072748322082193617
$ws.Range("G5").Value = -0.064571424820272477
$ws.Range("H5").Value = 1.374273709118129

# row 6
$ws.Range("B6").Value = 0.00021757471112746261
$ws.Range("C6").Value = 0.0148037399288559
$ws.Range("D6").Value = 0.00038138826022926509
$ws.Range("E6").Value = -0.074192312159272689
$ws.Range("F6").Value = 0.1200715718075074
$ws.Range("G6").Value = 0.0140742384611843
$ws.Range("H6").Value = 1.578016206147302

# row 7
$ws.Range("B7").Value = 0.00060236784881574656
$ws.Range("C7").Value = 0.014203065137333619
$ws.Range("D7").Value = 0.0012405417349583561
$ws.Range("E7").Value = -0.10453832109995891
$ws.Range("F7").Value = 0.084002952088621186
$ws.Range("G7").Value = -0.034522049307325893
$ws.Range("H7").Value = 1.354151788073279

# row 8
$ws.Range("B8").Value = -0.0052913224801868358
$ws.Range("C8").Value = 0.1718285009853856
$ws.Range("D8").Value = 0.002202124972804143
$ws.Range("E8").Value = -4.5926297537767873
$ws.Range("F8").Value = 0.1081991258551796
$ws.Range("G8").Value = -0.076241159022769373
$ws.Range("H8").Value = 1.4696348066607869

# row 9
$ws.Range("B9").Value = -0.00016170070126910729
$ws.Range("C9").Value = 0.01411320038918913
$ws.Range("D9").Value = -0.000062108568302576828
$ws.Range("E9").Value = -0.071814307724510584
$ws.Range("F9").Value = 0.075481880201765605
$ws.Range("G9").Value = -0.02352010970556179
$ws.Range("H9").Value = 1.498410961279222

# row 10
$ws.Range("B10").Value = 0.0019333823778251509
$ws.Range("C10").Value = 0.062170271363395438
$ws.Range("D10").Value = 0.0019276925523570301
$ws.Range("E10").Value = -0.46473017535485012
$ws.Range("F10").Value = 0.27063751771878591
$ws.Range("G10").Value = 0.036385275176323971
$ws.Range("H10").Value = 1.8833039870540731

# row 11
$ws.Range("B11").Value = 0.00021726506607250309
$ws.Range("C11").Value = 0.016122407180638371
$ws.Range("D11").Value = 0.00071124117500787065
$ws.Range("E11").Value = -0.1309834911294932
$ws.Range("F11").Value = 0.16421456444571089
$ws.Range("G11").Value = 0.027490006976122389
$ws.Range("H11").Value = 1.466031109924135

# row 12
$ws.Range("B12").Value = 0.000078995316675032482
$ws.Range("C12").Value = 0.017287155481237441
$ws.Range("D12").Value = 0.0010452801560383309
$ws.Range("E12").Value = -0.1854610663248373
$ws.Range("F12").Value = 0.1159782860158725
$ws.Range("G12").Value = -0.07801764312247994
$ws.Range("H12").Value = 1.5610892419102029

# row 13
$ws.Range("B13").Value = -0.000043159533537880947
$ws.Range("C13").Value = 0.013446727882649929
$ws.Range("D13").Value = 0.00027145805224026992
$ws.Range("E13").Value = -0.1151170613425148
$ws.Range("F13").Value = 0.12125382552797249
$ws.Range("G13").Value = 0.01522940082392734
$ws.Range("H13").Value = 1.3671551798247259

# row 14
$ws.Range("B14").Value = 0.001244029812180135
$ws.Range("C14").Value = 0.047346124650061883
$ws.Range("D14").Value = 0.0029006985243920091
$ws.Range("E14").Value = -0.28179629930639521
$ws.Range("F14").Value = 0.7367530178235584
$ws.Range("G14").Value = -0.096659339330555039
$ws.Range("H14").Value = 1.5850787190258671

# row 15
$ws.Range("B15").Value = 0.00028711559839578491
$ws.Range("C15").Value = 0.013110437925915
$ws.Range("D15").Value = 0.00075768555327559994
$ws.Range("E15").Value = -0.1317579985597099
$ws.Range("F15").Value = 0.094088164058783619
$ws.Range("G15").Value = -0.0047319044193624453
$ws.Range("H15").Value = 1.4606284524625079

# row 16
$ws.Range("B16").Value = 0.000084839111695335766
$ws.Range("C16").Value = 0.015802723927694742
$ws.Range("D16").Value = 0.00057513513271789662
$ws.Range("E16").Value = -0.13054858703866051
$ws.Range("F16").Value = 0.14970686160126689
$ws.Range("G16").Value = 0.0018067848035785931
$ws.Range("H16").Value = 1.5779994320982209

# row 17
$ws.Range("B17").Value = 0.0015138978236896891
$ws.Range("C17").Value = 0.035308498205441359
$ws.Range("D17").Value = 0.001868980348497473
$ws.Range("E17").Value = -0.37756612018803892
$ws.Range("F17").Value = 0.16851379041827741
$ws.Range("G17").Value = -0.011808629433778721
$ws.Range("H17").Value = 1.4106050605859359

# row 18
$ws.Range("B18").Value = 0.0012194034131696109
$ws.Range("C18").Value = 0.021265379112407411
$ws.Range("D18").Value = 0.001922094626221504
$ws.Range("E18").Value = -0.2057265653202078
$ws.Range("F18").Value = 0.1339103958505872
$ws.Range("G18").Value = -0.032563139561844638
$ws.Range("H18").Value = 1.3685283859815549

# row 19
$ws.Range("B19").Value = 0.0002719931329485087
$ws.Range("C19").Value = 0.01491592731914511
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = -0.2009204145294996
$ws.Range("F19").Value = 0.047447093372394278
$ws.Range("G19").Value = 0.16154558992285559
$ws.Range("H19").Value = 3.2688256559461499

# row 20
$ws.Range("B20").Value = 0.0014494612591920339
$ws.Range("C20").Value = 0.024961353148912419
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = -0.092255489394731605
$ws.Range("F20").Value = 0.2252276481622868
$ws.Range("G20").Value = 0.049451987481936847
$ws.Range("H20").Value = 1.405111042272472

# row 21
$ws.Range("B21").Value = 0.0004510802611599478
$ws.Range("C21").Value = 0.01181331735697087
$ws.Range("D21").Value = 0.00053654856223106506
$ws.Range("E21").Value = -0.051069383901130472
$ws.Range("F21").Value = 0.1097401996827569
$ws.Range("G21").Value = 0.052835926696156582
$ws.Range("H21").Value = 1.5099550271627471

# row 22
$ws.Range("B22").Value = 0.0024151700762201889
$ws.Range("C22").Value = 0.079071529143920669
$ws.Range("D22").Value = 0.001492403001892129
$ws.Range("E22").Value = -0.55073174413121784
$ws.Range("F22").Value = 0.32497054552207327
$ws.Range("G22").Value = 0.076682003214706071
$ws.Range("H22").Value = 1.653084038250761

# row 23
$ws.Range("B23").Value = 0.000013058730819595251
$ws.Range("C23").Value = 0.01460929575752308
$ws.Range("D23").Value = 0.00078666551698924891
$ws.Range("E23").Value = -0.087669717954606874
$ws.Range("F23").Value = 0.14568249246448681
$ws.Range("G23").Value = -0.052966219299803191
$ws.Range("H23").Value = 1.5336515701112849

# row 24
$ws.Range("B24").Value = 0.00038504117995412932
$ws.Range("C24").Value = 0.022365801416451961
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = -0.075670361506415151
$ws.Range("F24").Value = 0.20961677501600781
$ws.Range("G24").Value = -0.013784849987351159
$ws.Range("H24").Value = 1.2646956442320809

# row 25
$ws.Range("B25").Value = 0.000037589637523131212
$ws.Range("C25").Value = 0.01301362000735776
$ws.Range("D25").Value = -0.0001568137249838841
$ws.Range("E25").Value = -0.067724656807788719
$ws.Range("F25").Value = 0.053358431674180203
$ws.Range("G25").Value = 0.043656162520373937
$ws.Range("H25").Value = 1.357568990205434

# row 26
$ws.Range("B26").Value = -0.0000065356340833683287
$ws.Range("C26").Value = 0.0054508079336273423
$ws.Range("D26").Value = -0.000045497760382377053
$ws.Range("E26").Value = -0.052569703281064813
$ws.Range("F26").Value = 0.053393347829450218
$ws.Range("G26").Value = -0.042615520338345383
$ws.Range("H26").Value = 2.4380424533751128

# row 27
$ws.Range("B27").Value = 0.000049670859439528568
$ws.Range("C27").Value = 0.023453049587518589
$ws.Range("D27").Value = 0.001060309450344743
$ws.Range("E27").Value = -0.40467437019784919
$ws.Range("F27").Value = 0.18261945405023641
$ws.Range("G27").Value = -0.035084951284693362
$ws.Range("H27").Value = 1.31514876552308

# row 28
$ws.Range("B28").Value = 0.00061548055460312011
$ws.Range("C28").Value = 0.014787046361160971
$ws.Range("D28").Value = 0.0011478884345348379
$ws.Range("E28").Value = -0.099944852300803078
$ws.Range("F28").Value = 0.088808406952074925
$ws.Range("G28").Value = 0.029557928680021821
$ws.Range("H28").Value = 1.610515751422682

# --- View / formatting changes ---
# ColumnWidth is pixel-quantized by Excel; 14.6 is the closest input that
# snaps to the target rendered width (~15.57 chars) for column A.
$ws.Columns("A").ColumnWidth = 14.6
$ws.Range("J4").Select()
